# Insert a new data row at row 747 (pushing the existing rows 747:788 down
# to 748:789) and populate it with the new record:
#   A747 = "2026/01/31"  (text, same formatting as the date-text column)
#   B747 = "土"
#   C747 = 19
#   D747 = 201
#
# This mirrors the source edit, which inserted one additional reading
# (2026/01/31, 土, 19, 201) just before the "2026/12/29" block, shifting
# every subsequent row down by one and growing the sheet dimension from
# A1:D788 to A1:D789.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 747:788 down to 748:789, leaving a blank row 747 behind.
$ws.Rows.Item(747).Insert()

# Column A in this sheet stores the date as literal text (e.g. "2026/12/29"),
# not a real Excel date value. Assigning a date-shaped string straight to
# .Value would make Excel auto-convert it into a date serial number, so we
# instead build it as a text formula, then copy/paste the resulting value
# back over itself (values only) to end up with a plain text cell that
# carries no formula and no extra/explicit formatting.
$ws.Range("A747").Formula = "=""2026/01/31"""
$ws.Range("A747").Copy()
$ws.Range("A747").PasteSpecial(-4163)

$ws.Range("B747").Value = "土"
$ws.Range("C747").Value = 19
$ws.Range("D747").Value = 201
